# Update MSE (B) and R2 (C) columns with values computed using the
# Mean Absolute Error metric instead of the previous metric.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ B = 0.5383957249147062;  C = 0.9892793977683405 }
    3  = @{ B = 0.2471851971389369;  C = 0.9951687591056781 }
    4  = @{ B = 0.2793909772648969;  C = 0.9946251388627957 }
    5  = @{ B = 0.4141130016185339;  C = 0.9918342585793326 }
    6  = @{ B = 0.5028005785757672;  C = 0.9852314298888494 }
    7  = @{ B = 0.09478778173813857; C = 0.9986900071912249 }
    8  = @{ B = 0.03292551676276323; C = 0.9996606437193195 }
    9  = @{ B = 0.09714035244815278; C = 0.9994213626368167 }
    10 = @{ B = 0.06788798627532909; C = 0.9987598472634609 }
    11 = @{ B = 0.1229159912292265;  C = 0.9909185187883635 }
    12 = @{ B = 0.05125156037683928; C = 0.9984927790469791 }
    13 = @{ B = 0.06048728111847583; C = 0.9994262510313371 }
    14 = @{ B = 0.05721954093985681; C = 0.9992252827361765 }
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row].B
    $ws.Range("C$row").Value = $updates[$row].C
}
